{"js": "// Edit of SDD regarding MVP: a Presenter only controls a single View (like a\n// regular Controller), but only passes communicative information. Applies the\n// full set of textual edits described by the commit diff.\n\nconst doc = context.document;\nconst body = doc.body;\n\nasync function replaceOnce(searchText, replaceText) {\n  const results = body.search(searchText, { matchCase: true, ignoreSpace: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + searchText);\n  }\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Cached DATE field result: Thursday 27 April 2017 -> Saturday 29 April 2017 (Swedish).\nawait replaceOnce(\n  \"torsdag den 27 april 2017\",\n  \"l\u00f6rdag den 29 april 2017\"\n);\n\n// 2. \"design pattern\" -> \"design principle\"\nawait replaceOnce(\n  \"design pattern. In an Android application, all MVC-based architectural flows will look similar to this \",\n  \"design principle. In an Android application, all MVC-based architectural flows will look similar to this \"\n);\n\n// 3. Drop the word \"instead\" before \"only handles\".\nawait replaceOnce(\n  \"the data of the Model, the Presenter instead only handles \",\n  \"the data of the Model, the Presenter only handles \"\n);\n\n// 4. Rewrite the sentence about what the Presenter keeps track of / the View displaying itself.\nawait replaceOnce(\n  \" More precisely, the Presenter keeps track of what View to display to the User, but the View itself handles the logic for how it should display itself.\",\n  \" More precisely, the Presenter keeps track of the user interaction in the View and modifies the Model accordingly, but the View itself handles the logic for what it should display.\"\n);\n\n// 5. Rewrite how data flows from the Model to the Activity/View.\nawait replaceOnce(\n  \"rom Model by the Presenter and handed to the specific Activity for the View, which and then inflates the necessary \",\n  \"rom Model by the Presenter and handed to its associated Activity in the View. The Activity then inflates the necessary \"\n);\n\n// 6. Replace description of the Presenter's role relative to the Activity/Controller.\nawait replaceOnce(\n  \"An Activity sort of acts as a small Controller for each specific View, and the Presenter is the managing operator which tells each Controller when to do their job. This way, the V\",\n  \"An Activity sort of acts as a small Controller for each specific View, and the Presenter is the communicator between the Activity and the Model, only telling information that is necessary to be shared between the two. This way, the V\"\n);\n\n// 7. Insert \"for each View \" so a single Presenter is scoped to one View, like a Controller.\nawait replaceOnce(\n  \"letting the Controller be the Presenter in the form of\",\n  \"letting the Controller for each View be the Presenter in the form of\"\n);\n\n// 7b. Word's \"_GoBack\" bookmark tracks the location of the last edit. Since the\n// last edit of that paragraph was typing \"for each View \", move the bookmark\n// there (right before \"be the Presenter in the form of\"), matching real Word\n// behaviour, instead of leaving it at the paragraph's previous edit spot.\ndoc.deleteBookmark(\"_GoBack\");\nconst gobackTarget = body.search(\"be the Presenter in the form of\", { matchCase: true });\ngobackTarget.load(\"items\");\nawait context.sync();\nif (gobackTarget.items.length > 0) {\n  const gobackRange = gobackTarget.items[0].getRange(\"Start\");\n  gobackRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Edit of SDD regarding MVP. A presenter is only controlling one single view,\n# such as a regular Controller, but only communicative information.\n#\n# Applies the full set of textual edits described by the commit diff using\n# the Word COM object model (Find/Replace + Bookmarks).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $ok = $find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $FindText\"\n    }\n}\n\n# 1. Cached DATE field result: Thursday 27 April 2017 -> Saturday 29 April 2017 (Swedish).\nReplace-Text \"torsdag den 27 april 2017\" \"l\u00f6rdag den 29 april 2017\"\n\n# 2. \"design pattern\" -> \"design principle\"\nReplace-Text \"design pattern. In an Android application, all MVC-based architectural flows will look similar to this \" \"design principle. In an Android application, all MVC-based architectural flows will look similar to this \"\n\n# 3. Drop the word \"instead\" before \"only handles\".\nReplace-Text \"the data of the Model, the Presenter instead only handles \" \"the data of the Model, the Presenter only handles \"\n\n# 4. Rewrite the sentence about what the Presenter keeps track of / the View displaying itself.\nReplace-Text \" More precisely, the Presenter keeps track of what View to display to the User, but the View itself handles the logic for how it should display itself.\" \" More precisely, the Presenter keeps track of the user interaction in the View and modifies the Model accordingly, but the View itself handles the logic for what it should display.\"\n\n# 5. Rewrite how data flows from the Model to the Activity/View.\nReplace-Text \"rom Model by the Presenter and handed to the specific Activity for the View, which and then inflates the necessary \" \"rom Model by the Presenter and handed to its associated Activity in the View. The Activity then inflates the necessary \"\n\n# 6. Replace description of the Presenter's role relative to the Activity/Controller.\nReplace-Text \"An Activity sort of acts as a small Controller for each specific View, and the Presenter is the managing operator which tells each Controller when to do their job. This way, the V\" \"An Activity sort of acts as a small Controller for each specific View, and the Presenter is the communicator between the Activity and the Model, only telling information that is necessary to be shared between the two. This way, the V\"\n\n# 7. Insert \"for each View \" so a single Presenter is scoped to one View, like a Controller.\nReplace-Text \"letting the Controller be the Presenter in the form of\" \"letting the Controller for each View be the Presenter in the form of\"\n\n# 7b. Word's \"_GoBack\" bookmark tracks the location of the last edit. Since the\n# last edit of that paragraph was typing \"for each View \", move the bookmark\n# there (right before \"be the Presenter in the form of\"), matching real Word\n# behaviour, instead of leaving it at the paragraph's previous edit spot.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$findGoBack = $d.Content.Find\n$findGoBack.ClearFormatting()\n$foundGoBack = $findGoBack.Execute(\"be the Presenter in the form of\", $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($foundGoBack) {\n    $goBackTarget = $findGoBack.Parent\n    $goBackRange = $d.Range($goBackTarget.Start, $goBackTarget.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $goBackRange)\n}\n"}
